$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename advertising line items: remove Radio/Television/Direct Mail,
# add Banner/Ad-Words/Social in their place (rows 7-10), apply the
# updated "Comma0"-derived style to match the source look.
$ws.Range("B7").Value = "Banner"
$ws.Range("B8").Value = "Print"
$ws.Range("B10").Value = "Social"
$ws.Range("B9").Value = "Ad-Words"

$ws.Range("B7:B10").Style = "Comma0"

# Increase Videos budget figures across all four quarters
$ws.Range("C29").Value = 10000
$ws.Range("E29").Value = 10000
$ws.Range("G29").Value = 10000
$ws.Range("I29").Value = 10000

# Move the active selection to B51, as recorded in the saved view state
$ws.Range("B51").Select()
